$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("H41").Value = 521
$ws.Range("I41").Value = 486.7143
$ws.Range("J41").Value = 564.63635
$ws.Range("K41").Value = 486.7143
$ws.Range("L41").Value = 564.63635
$ws.Range("M41").Value = -46.71429999999998
$ws.Range("N41").Value = -1444.63635
$ws.Range("H76").Value = 14294714
$ws.Range("I76").Value = 16674666
$ws.Range("J76").Value = 15000
$ws.Range("K76").Value = 16674666
$ws.Range("L76").Value = 15000
$ws.Range("M76").Value = -16674351
$ws.Range("N76").Value = -15630
$ws.Range("H79").Value = 14294714
$ws.Range("I79").Value = 16674666
$ws.Range("J79").Value = 15000
$ws.Range("K79").Value = 16674666
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = -16673574
$ws.Range("N79").Value = -17184
$ws.Range("H86").Value = 9927.571
$ws.Range("I86").Value = 9698.6
$ws.Range("K86").Value = 9698.6
$ws.Range("M86").Value = -8575.6
$ws.Range("H89").Value = 9927.571
$ws.Range("I89").Value = 9698.6
$ws.Range("K89").Value = 48493
$ws.Range("M89").Value = -42877
$ws.Range("H113").Value = 251751.5
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 501003
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 501003
$ws.Range("M113").Value = 754
$ws.Range("N113").Value = -507511
$ws.Range("H116").Value = 5198.25
$ws.Range("I116").Value = 3996.5
$ws.Range("K116").Value = 3996.5
$ws.Range("M116").Value = -554.5
$ws.Range("H138").Value = 4966.328
$ws.Range("J138").Value = 6164.619
$ws.Range("L138").Value = 18493.857
$ws.Range("N138").Value = -28773.857

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23398.328
$ws.Range("I32").Value = 23853.518
$ws.Range("J32").Value = 19496.715
$ws.Range("K32").Value = 23853.518
$ws.Range("L32").Value = 19496.715
$ws.Range("M32").Value = -23566.518
$ws.Range("N32").Value = -20070.715
$ws.Range("H74").Value = 273196.06
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31748
$ws.Range("H77").Value = 273196.06
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -158736

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 83335610
$ws.Range("I105").Value = 125002270
$ws.Range("J105").Value = 2269.75
$ws.Range("K105").Value = 125002270
$ws.Range("L105").Value = 2269.75
$ws.Range("M105").Value = -125000523
$ws.Range("N105").Value = -5763.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25647618
$ws.Range("I31").Value = 142860380
$ws.Range("K31").Value = 142860380
$ws.Range("M31").Value = -142860085
$ws.Range("H34").Value = 25647618
$ws.Range("I34").Value = 142860380
$ws.Range("K34").Value = 142860380
$ws.Range("M34").Value = -142860178
$ws.Range("H62").Value = 28333.334
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 28333.334
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 28333.334
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -29581.334
$ws.Range("H65").Value = 28333.334
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 28333.334
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 141666.67
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -147906.67
$ws.Range("H99").Value = 2791.3157
$ws.Range("I99").Value = 2753.7
$ws.Range("J99").Value = 2833.111
$ws.Range("K99").Value = 2753.7
$ws.Range("L99").Value = 2833.111
$ws.Range("M99").Value = -1255.7
$ws.Range("N99").Value = -5829.111
$ws.Range("H126").Value = 2791.3157
$ws.Range("I126").Value = 2753.7
$ws.Range("J126").Value = 2833.111
$ws.Range("K126").Value = 8261.099999999999
$ws.Range("L126").Value = 8499.332999999999
$ws.Range("M126").Value = -5791.099999999999
$ws.Range("N126").Value = -13439.333
$ws.Range("H132").Value = 28809.809
$ws.Range("I132").Value = 1860.5333
$ws.Range("K132").Value = 5581.5999
$ws.Range("M132").Value = -3051.5999
$ws.Range("H141").Value = 159104.92
$ws.Range("J141").Value = 181764.42
$ws.Range("L141").Value = 181764.42
$ws.Range("N141").Value = -192124.42

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 50004496
$ws.Range("J122").Value = 1027.4117
$ws.Range("L122").Value = 9246.705300000001
$ws.Range("N122").Value = -14146.7053

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7475
$ws.Range("I80").Value = 5962.5
$ws.Range("K80").Value = 5962.5
$ws.Range("M80").Value = -4964.5
$ws.Range("H83").Value = 7475
$ws.Range("I83").Value = 5962.5
$ws.Range("K83").Value = 29812.5
$ws.Range("M83").Value = -24820.5
$ws.Range("H102").Value = 1609.7
$ws.Range("I102").Value = 1188
$ws.Range("K102").Value = 1188
$ws.Range("M102").Value = 434
$ws.Range("H122").Value = 1263.5454
$ws.Range("I122").Value = 1058.5555
$ws.Range("K122").Value = 3175.6665
$ws.Range("M122").Value = -725.6664999999998

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1499
$ws.Range("I7").Value = 1499
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1499
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1387
$ws.Range("N7").ClearContents()
$ws.Range("H16").Value = 1519.5
$ws.Range("I16").Value = 1052.35
$ws.Range("J16").Value = 3855.25
$ws.Range("K16").Value = 1052.35
$ws.Range("L16").Value = 3855.25
$ws.Range("M16").Value = -882.3499999999999
$ws.Range("N16").Value = -4195.25
$ws.Range("H40").Value = 8859.532999999999
$ws.Range("I40").Value = 10322.333
$ws.Range("K40").Value = 10322.333
$ws.Range("M40").Value = -10186.333
$ws.Range("H46").Value = 7429.1875
$ws.Range("I46").Value = 1482.6666
$ws.Range("J46").Value = 10997.1
$ws.Range("K46").Value = 1482.6666
$ws.Range("L46").Value = 10997.1
$ws.Range("M46").Value = -1294.6666
$ws.Range("N46").Value = -11373.1
$ws.Range("H61").Value = 6166.625
$ws.Range("I61").Value = 6333.357
$ws.Range("J61").Value = 4999.5
$ws.Range("K61").Value = 6333.357
$ws.Range("L61").Value = 4999.5
$ws.Range("M61").Value = -6131.357
$ws.Range("N61").Value = -5403.5
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50450
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51560
$ws.Range("H82").Value = 2195.6667
$ws.Range("I82").Value = 2050
$ws.Range("J82").Value = 2487
$ws.Range("K82").Value = 2050
$ws.Range("L82").Value = 2487
$ws.Range("M82").Value = -1689
$ws.Range("N82").Value = -3209
$ws.Range("H85").Value = 2195.6667
$ws.Range("I85").Value = 2050
$ws.Range("J85").Value = 2487
$ws.Range("K85").Value = 2050
$ws.Range("L85").Value = 2487
$ws.Range("M85").Value = -802
$ws.Range("N85").Value = -4983
$ws.Range("H113").Value = 6166.625
$ws.Range("I113").Value = 6333.357
$ws.Range("J113").Value = 4999.5
$ws.Range("K113").Value = 6333.357
$ws.Range("L113").Value = 4999.5
$ws.Range("M113").Value = -4163.357
$ws.Range("N113").Value = -9339.5
$ws.Range("H122").Value = 62507268
$ws.Range("I122").Value = 71435750
$ws.Range("K122").Value = 214307250
$ws.Range("M122").Value = -214304800
$ws.Range("H126").Value = 1499
$ws.Range("I126").Value = 1499
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4497
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2027
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 6282.021
$ws.Range("I136").Value = 4672.6284
$ws.Range("K136").Value = 14017.8852
$ws.Range("M136").Value = -11467.8852

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12273.1
$ws.Range("I81").Value = 5311.636
$ws.Range("J81").Value = 14913.655
$ws.Range("K81").Value = 10623.272
$ws.Range("L81").Value = 29827.31
$ws.Range("M81").Value = -9562.272000000001
$ws.Range("N81").Value = -31949.31
$ws.Range("H84").Value = 12273.1
$ws.Range("I84").Value = 5311.636
$ws.Range("J84").Value = 14913.655
$ws.Range("K84").Value = 53116.36
$ws.Range("L84").Value = 149136.55
$ws.Range("M84").Value = -47812.36
$ws.Range("N84").Value = -159744.55
$ws.Range("H122").Value = 3483.205
$ws.Range("I122").Value = 3772.4075
$ws.Range("K122").Value = 11317.2225
$ws.Range("M122").Value = -8867.2225
$ws.Range("H126").Value = 1710.5
$ws.Range("I126").Value = 1815.7
$ws.Range("J126").Value = 1184.5
$ws.Range("K126").Value = 5447.1
$ws.Range("L126").Value = 3553.5
$ws.Range("M126").Value = -2977.1
$ws.Range("N126").Value = -8493.5
